$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in column D (sex = "F") and column E (environ = 1) for rows 2-10
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 4).Value = "F"
    $ws.Cells.Item($r, 5).Value = 1
}

# Update selection to match the new active range
$ws.Range("D2:E10").Select()
